$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in zeros for previously-empty cure/death/new-case cells so every
# region row has a value (supports drawing a bar chart for all regions).
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0

$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 0

$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0

$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 0

$ws.Range("B7").Value = 0
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 0

$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 0

$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 0

$ws.Range("C10").Value = 0
$ws.Range("D10").Value = 0

$ws.Range("D11").Value = 0

$ws.Range("B12").Value = 0
$ws.Range("C12").Value = 0
$ws.Range("D12").Value = 0

# Move the marker/selection to C15:C16 (active cell C15), matching the
# "add marker in line" part of the edit.
$ws.Range("C15:C16").Select()
